$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: Rewrite the summary paragraph ("Computer Science graduate...")
# -----------------------------------------------------------------
$rng = $d.Content
$oldSummary = "Computer Science graduate, experienced with object-oriented programming and full-stack development. Pursuing full-time opportunities in the industry."
$null = $rng.Find.Execute($oldSummary)
$pStart = $rng.Start

$nbHyphen = [char]0x2011
$newSummary = "Detail-oriented summa cum laude Computer Science graduate with full" + $nbHyphen + `
    "stack and mobile development experience in C#, JS/React, Dart/Flutter, and cloud platforms (Oracle, GCP, AWS). Delivered production" + $nbHyphen + `
    "ready solutions like high" + $nbHyphen + `
    "traffic sites, CI/CD systems, a NASA" + $nbHyphen + `
    "sponsored game, and scalable APIs. Eager to leverage problem" + $nbHyphen + `
    "solving and agile collaboration skills."

$rng.Text = $newSummary

# Italicize "summa cum laude"
$italicRange = $d.Range($pStart + 16, $pStart + 31)
$italicRange.Font.Italic = 1

# Apply the Cambria Math font to the 5 non-breaking hyphens (matches the
# special-character insertion Word performs for this glyph)
$hyphenOffsets = @(67, 199, 225, 262, 323)
foreach ($off in $hyphenOffsets) {
    $hr = $d.Range($pStart + $off, $pStart + $off + 1)
    $hr.Font.Name = "Cambria Math"
}

# -----------------------------------------------------------------
# Change 2: "Web Development Intern" -> "Web Developer"
# -----------------------------------------------------------------
$rng2 = $d.Content
$null = $rng2.Find.Execute("Web Development Intern")
$wStart = $rng2.Start
# "Web Develop" is 11 chars; replace the remaining "ment Intern" (11 chars) with "er"
$tailRange = $d.Range($wStart + 11, $wStart + 22)
$tailRange.Text = "er"

# -----------------------------------------------------------------
# Change 3: mark a rendered page break right before the "Practic(ed)" bullet
# -----------------------------------------------------------------
$rng3 = $d.Content
$null = $rng3.Find.Execute("Practiced Agile")
$practicStart = $rng3.Start
$breakPoint = $d.Range($practicStart, $practicStart)
$breakXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$breakPoint.InsertXML($breakXml)
